# The customer-facing "Total" column on the summary sheet is being
# relabeled to clarify that the figure is the per-room total (this
# workbook is emailed straight to the customer as their results file).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Total Per Room"

# Cosmetic leftovers from the editing session (widen the label column so
# the new, longer header fits; leave the cursor where the author left it).
$ws.Columns.Item(1).ColumnWidth = 31.42
$ws.Range("G9").Select() | Out-Null
